$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 9.423641
$ws.Range("H2").Value = 28.270923
$ws.Range("I2").Value = 0.3411090217977475
$ws.Range("J2").Value = 0.3411090217977475
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 14.96835266666667
$ws.Range("N2").Value = 44.905058
$ws.Range("O2").Value = 0.1240053612000741
$ws.Range("P2").Value = 0.1240053612000741
$ws.Range("Q2").Value = 141.0563818920593
$ws.Range("R2").Value = 1269.507437028534
$ws.Range("S2").Value = 0.04229934745663363
$ws.Range("T2").Value = 0.04229934745663363

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 9.423641
$ws.Range("H3").Value = 28.270923
$ws.Range("I3").Value = 0.3411090217977475
$ws.Range("J3").Value = 0.3411090217977475
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 80.77474466666666
$ws.Range("N3").Value = 242.324234
$ws.Range("O3").Value = 0.6691786071115035
$ws.Range("P3").Value = 0.6691786071115035
$ws.Range("Q3").Value = 761.1921956053313
$ws.Range("R3").Value = 6850.729760447982
$ws.Range("S3").Value = 0.2282628600797842
$ws.Range("T3").Value = 0.2282628600797842

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 9.423641
$ws.Range("H4").Value = 28.270923
$ws.Range("I4").Value = 0.3411090217977475
$ws.Range("J4").Value = 0.3411090217977475
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 24.96420533333334
$ws.Range("N4").Value = 74.892616
$ws.Range("O4").Value = 0.2068160316884225
$ws.Range("P4").Value = 0.2068160316884225
$ws.Range("Q4").Value = 235.2537089116187
$ws.Range("R4").Value = 2117.283380204568
$ws.Range("S4").Value = 0.07054681426132975
$ws.Range("T4").Value = 0.07054681426132976

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 14.13955433333333
$ws.Range("H5").Value = 42.418663
$ws.Range("I5").Value = 0.5118116816312757
$ws.Range("J5").Value = 0.5118116816312757
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 14.96835266666667
$ws.Range("N5").Value = 44.905058
$ws.Range("O5").Value = 0.1240053612000741
$ws.Range("P5").Value = 0.1240053612000741
$ws.Range("Q5").Value = 211.6458358108282
$ws.Range("R5").Value = 1904.812522297454
$ws.Range("S5").Value = 0.06346739244710366
$ws.Range("T5").Value = 0.06346739244710366

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 14.13955433333333
$ws.Range("H6").Value = 42.418663
$ws.Range("I6").Value = 0.5118116816312757
$ws.Range("J6").Value = 0.5118116816312757
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 80.77474466666666
$ws.Range("N6").Value = 242.324234
$ws.Range("O6").Value = 0.6691786071115035
$ws.Range("P6").Value = 0.6691786071115035
$ws.Range("Q6").Value = 1142.11889097546
$ws.Range("R6").Value = 10279.07001877914
$ws.Range("S6").Value = 0.3424934282174133
$ws.Range("T6").Value = 0.3424934282174133

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 14.13955433333333
$ws.Range("H7").Value = 42.418663
$ws.Range("I7").Value = 0.5118116816312757
$ws.Range("J7").Value = 0.5118116816312757
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 24.96420533333334
$ws.Range("N7").Value = 74.892616
$ws.Range("O7").Value = 0.2068160316884225
$ws.Range("P7").Value = 0.2068160316884225
$ws.Range("Q7").Value = 352.9827376991565
$ws.Range("R7").Value = 3176.844639292408
$ws.Range("S7").Value = 0.1058508609667587
$ws.Range("T7").Value = 0.1058508609667587

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.063282999999999
$ws.Range("H8").Value = 12.189849
$ws.Range("I8").Value = 0.1470792965709768
$ws.Range("J8").Value = 0.1470792965709768
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 14.96835266666667
$ws.Range("N8").Value = 44.905058
$ws.Range("O8").Value = 0.1240053612000741
$ws.Range("P8").Value = 0.1240053612000741
$ws.Range("Q8").Value = 60.82065292847132
$ws.Range("R8").Value = 547.3858763562419
$ws.Range("S8").Value = 0.0182386212963368
$ws.Range("T8").Value = 0.0182386212963368

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.063282999999999
$ws.Range("H9").Value = 12.189849
$ws.Range("I9").Value = 0.1470792965709768
$ws.Range("J9").Value = 0.1470792965709768
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 80.77474466666666
$ws.Range("N9").Value = 242.324234
$ws.Range("O9").Value = 0.6691786071115035
$ws.Range("P9").Value = 0.6691786071115035
$ws.Range("Q9").Value = 328.2106468334073
$ws.Range("R9").Value = 2953.895821500666
$ws.Range("S9").Value = 0.09842231881430601
$ws.Range("T9").Value = 0.09842231881430601

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.063282999999999
$ws.Range("H10").Value = 12.189849
$ws.Range("I10").Value = 0.1470792965709768
$ws.Range("J10").Value = 0.1470792965709768
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 24.96420533333334
$ws.Range("N10").Value = 74.892616
$ws.Range("O10").Value = 0.2068160316884225
$ws.Range("P10").Value = 0.2068160316884225
$ws.Range("Q10").Value = 101.4366311394427
$ws.Range("R10").Value = 912.929680254984
$ws.Range("S10").Value = 0.03041835646033404
$ws.Range("T10").Value = 0.03041835646033404
